$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param(
        [int]$RowNumber,
        [object[]]$Values
    )

    for ($i = 0; $i -lt $Values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($RowNumber, $col)

        if ($col -eq 2) {
            # "Date" column: force text so Excel does not auto-convert the
            # dd/mm/yyyy looking string into a date serial number.
            $cell.NumberFormat = "@"
            $cell.Value = $Values[$i]
            $cell.ClearFormats()
        }
        else {
            $cell.Value = $Values[$i]
        }
    }
}

# ---------------------------------------------------------------------------
# 1) Row 3: replace the existing match (Atlante - Tampico Madero) with the
#    new match (Bucaramanga - Ind. Medellin).
# ---------------------------------------------------------------------------
$row3Data = @(
    "2TUEylld","03/11/2024","22:30","COLOMBIA - PRIMERA A","Bucaramanga","Ind. Medellin",
    2.3,2.8,3.6,3.1,1.91,4.33,1.11,6.5,1.5,2.5,2.6,1.48,1.57,2.25,2.1,1.67,
    6,9.5,10,21,23,41,6,5.5,19,67,8,15,13,41,34,51,201,
    4,15,29,51,81,301,2.25,9,81,5,21,34,67,126,351,126,126
)
Set-RowData 3 $row3Data

# ---------------------------------------------------------------------------
# 2) Insert a brand-new row 6 (Vancouver Whitecaps - Los Angeles FC), pushing
#    the existing row 6 (New Mexico - Phoenix Rising) down to row 7.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Insert()

$row6Data = @(
    "MVdT00pJ","03/11/2024","22:45","USA - MLS","Vancouver Whitecaps","Los Angeles FC",
    2.45,3.7,2.7,3,2.3,3.2,1.03,15,1.2,4.33,1.67,2.15,1.33,3.25,1.54,2.25,
    10,13,9.5,23,17,23,15,7,13,41,11,15,10,29,19,23,151,
    4.75,13,19,41,51,126,3.25,7,41,5,15,21,41,51,126,401,151
)
Set-RowData 6 $row6Data
